$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Periodo Mora" (debt period) values for the four worker rows
# in the account-statement table: previous periods are dropped and the
# newest periods (2304-2306) are added, with "2303" moving down.
$ws.Range("E16").Value = "2306"
$ws.Range("E17").Value = "2305"
$ws.Range("E18").Value = "2304"
$ws.Range("E19").Value = "2303"
